$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = "R_SVA"
$ws.Range("B2").Value2 = [double]"2714888.40129829"
$ws.Range("B3").Value2 = [double]"2684356.106899888"
$ws.Range("A4").Value = "R_DYNHT"
$ws.Range("B4").Value2 = [double]"530519.5946675759"
$ws.Range("A5").Value = "R_O2Sat"
$ws.Range("B5").Value2 = [double]"396759.2753411022"
$ws.Range("A6").Value = "R_O2"
$ws.Range("B6").Value2 = [double]"323055.8579581335"
$ws.Range("A7").Value = "R_Depth"
$ws.Range("B7").Value2 = [double]"312594.9482566265"
$ws.Range("A8").Value = "R_PRES"
$ws.Range("B8").Value2 = [double]"310283.4061970197"
$ws.Range("A9").Value = "R_O2_sqrt"
$ws.Range("B9").Value2 = [double]"277008.4668065083"
$ws.Range("A10").Value = "R_PRES_cat_(317.5, inf]"
$ws.Range("B10").Value2 = [double]"249871.3379890129"
$ws.Range("A11").Value = "R_PO4"
$ws.Range("B11").Value2 = [double]"188769.7217200231"
$ws.Range("A12").Value = "R_PRES_cat_(0.0, 47.5]"
$ws.Range("B12").Value2 = [double]"186582.6372257432"
$ws.Range("A13").Value = "R_NO3"
$ws.Range("B13").Value2 = [double]"135301.1896482321"
$ws.Range("A14").Value = "R_Depth_cat_(315.5, 671.5]"
$ws.Range("B14").Value2 = [double]"114797.8075130635"
$ws.Range("A15").Value = "R_SALINITY"
$ws.Range("B15").Value2 = [double]"102144.5845694863"
$ws.Range("A16").Value = "R_SIO3"
$ws.Range("B16").Value2 = [double]"101280.6243001091"
$ws.Range("A17").Value = "R_PO4_cat_(0.0, 0.505]"
$ws.Range("B17").Value2 = [double]"98325.81600348902"
$ws.Range("A18").Value = "R_SIO3_cat_(0.0, 3.55]"
$ws.Range("B18").Value2 = [double]"77322.10248279784"
$ws.Range("A19").Value = "R_Depth_cat_(0.0, 15.5]"
$ws.Range("B19").Value2 = [double]"70912.69712999069"
$ws.Range("A20").Value = "R_Depth_cat_(671.5, inf]"
$ws.Range("B20").Value2 = [double]"59855.32231765103"
$ws.Range("A21").Value = "R_Depth_cat_(15.5, 47.5]"
$ws.Range("B21").Value2 = [double]"58327.66354584995"
$ws.Range("A22").Value = "R_CHLA_cat_(0.065, inf]"
$ws.Range("B22").Value2 = [double]"57300.55019632763"
$ws.Range("A23").Value = "R_PO4_cat_(2.505, inf]"
$ws.Range("B23").Value2 = [double]"54652.06555948611"
$ws.Range("A24").Value = "R_SIO3_cat_(42.05, inf]"
$ws.Range("B24").Value2 = [double]"44579.70349106604"
$ws.Range("A25").Value = "R_NO3_cat_(30.05, inf]"
$ws.Range("B25").Value2 = [double]"42233.56007720197"
$ws.Range("A26").Value = "R_CHLA_missing"
$ws.Range("B26").Value2 = [double]"29910.27770914006"
$ws.Range("A27").Value = "R_PHAEO_missing"
$ws.Range("B27").Value2 = [double]"29887.44598188028"
$ws.Range("A28").Value = "R_NO3_cat_(0.0, 0.55]"
$ws.Range("B28").Value2 = [double]"26783.15453731845"
$ws.Range("A29").Value = "R_PRES_cat_(118.5, 317.5]"
$ws.Range("B29").Value2 = [double]"26715.03066051535"
$ws.Range("A30").Value = "R_Depth_cat_(167.5, 315.5]"
$ws.Range("B30").Value2 = [double]"25822.77443068063"
$ws.Range("A31").Value = "Lat_Dec"
$ws.Range("B31").Value2 = [double]"24458.4687073659"
$ws.Range("A32").Value = "Phi"
$ws.Range("B32").Value2 = [double]"21039.15155552508"
$ws.Range("A33").Value = "R_PRES_cat_(47.5, 118.5]"
$ws.Range("B33").Value2 = [double]"16902.02769188509"
$ws.Range("A34").Value = "Rho"
$ws.Range("B34").Value2 = [double]"16370.32515075266"
$ws.Range("A35").Value = "R_Depth_cat_(47.5, 68.5]"
$ws.Range("B35").Value2 = [double]"13972.79318691749"
$ws.Range("A36").Value = "Lon_Dec"
$ws.Range("B36").Value2 = [double]"12528.24531616173"
$ws.Range("A37").Value = "R_PHAEO_cat_(0.015, 0.175]"
$ws.Range("B37").Value2 = [double]"12332.20428891879"
$ws.Range("A38").Value = "R_PHAEO_cat_(0.175, inf]"
$ws.Range("B38").Value2 = [double]"9823.414394576244"
$ws.Range("A39").Value = "R_PO4_cat_(1.025, 2.505]"
$ws.Range("B39").Value2 = [double]"8483.247546078766"
$ws.Range("A40").Value = "R_NO3_cat_(0.55, 10.85]"
$ws.Range("B40").Value2 = [double]"7144.489441342773"
$ws.Range("A41").Value = "R_SIO3_cat_(3.55, 10.05]"
$ws.Range("B41").Value2 = [double]"6592.698967501911"
$ws.Range("A42").Value = "R_PO4_cat_(0.505, 1.025]"
$ws.Range("B42").Value2 = [double]"6276.513287597355"
$ws.Range("A43").Value = "R_NO3_cat_(10.85, 30.05]"
$ws.Range("B43").Value2 = [double]"4969.737007009157"
$ws.Range("A44").Value = "R_NO2_cat_(0.045, 1.665]"
$ws.Range("B44").Value2 = [double]"4756.664337434716"
$ws.Range("A45").Value = "R_SIO3_cat_(10.05, 42.05]"
$ws.Range("B45").Value2 = [double]"4663.385373806259"
$ws.Range("A46").Value = "R_CHLA"
$ws.Range("B46").Value2 = [double]"3814.564890351276"
$ws.Range("A47").Value = "R_Depth_cat_(68.5, 117.5]"
$ws.Range("B47").Value2 = [double]"3546.943169601435"
$ws.Range("A48").Value = "R_PHAEO_cat_(0.005, 0.015]"
$ws.Range("B48").Value2 = [double]"1967.627482069903"
$ws.Range("A49").Value = "R_NH4_cat_(0.045, inf]"
$ws.Range("B49").Value2 = [double]"1591.697259690214"
$ws.Range("A50").Value = "R_PHAEO"
$ws.Range("B50").Value2 = [double]"1572.4400613119"
$ws.Range("A51").Value = "R_NO2"
$ws.Range("B51").Value2 = [double]"1551.306333968829"
$ws.Range("A52").Value = "R_Depth_cat_(117.5, 167.5]"
$ws.Range("B52").Value2 = [double]"1412.670400904832"
$ws.Range("C52").Value2 = [double]"1.422076229348544e-308"
$ws.Range("A53").Value = "R_NO2_missing"
$ws.Range("B53").Value2 = [double]"1397.864526552771"
$ws.Range("C53").Value2 = [double]"2.280935432565034e-305"
$ws.Range("A54").Value = "R_CHLA_cat_(0.0, 0.015]"
$ws.Range("B54").Value2 = [double]"1320.503532934729"
$ws.Range("C54").Value2 = [double]"1.281670447809801e-288"
$ws.Range("A55").Value = "R_SIO3_missing"
$ws.Range("B55").Value2 = [double]"1296.262630684214"
$ws.Range("C55").Value2 = [double]"2.275688171524229e-283"
$ws.Range("A56").Value = "R_PO4_missing"
$ws.Range("B56").Value2 = [double]"655.2476995483604"
$ws.Range("C56").Value2 = [double]"2.15686008779483e-144"
$ws.Range("A57").Value = "R_NO3_missing"
$ws.Range("B57").Value2 = [double]"569.8346465847553"
$ws.Range("C57").Value2 = [double]"7.596557236730098e-126"
$ws.Range("A58").Value = "R_NO2_cat_(0.005, 0.045]"
$ws.Range("B58").Value2 = [double]"555.84655807102"
$ws.Range("C58").Value2 = [double]"8.295775193515382e-123"
$ws.Range("B59").Value2 = [double]"469.3354099914836"
$ws.Range("C59").Value2 = [double]"5.189420311970072e-104"
$ws.Range("A60").Value = "R_NH4"
$ws.Range("B60").Value2 = [double]"464.6456392323434"
$ws.Range("C60").Value2 = [double]"5.424718510009815e-103"
$ws.Range("A61").Value = "R_NH4_cat_(0.005, 0.045]"
$ws.Range("B61").Value2 = [double]"155.721508817497"
$ws.Range("C61").Value2 = [double]"9.901815599470554e-36"
$ws.Range("A62").Value = "R_CHLA_cat_(0.055, 0.065]"
$ws.Range("B62").Value2 = [double]"152.0199182661288"
$ws.Range("C62").Value2 = [double]"6.372760923046343e-35"
$ws.Range("A63").Value = "R_NH4_missing"
$ws.Range("B63").Value2 = [double]"144.6229425417114"
$ws.Range("C63").Value2 = [double]"2.633928854950761e-33"
$ws.Range("A64").Value = "R_NO2_cat_(1.665, inf]"
$ws.Range("B64").Value2 = [double]"24.16611469822136"
$ws.Range("C64").Value2 = [double]"8.841122614622893e-07"
$ws.Range("A65").Value = "R_PHAEO_cat_(0.0, 0.005]"
$ws.Range("B65").Value2 = [double]"0"
$ws.Range("C65").Value2 = [double]"1"
$ws.Range("A66").Value = "R_NH4_cat_(0.0, 0.005]"
$ws.Range("A67").Value = "R_Depth_missing"
$ws.Range("A68").Value = "R_PRES_missing"
$ws.Range("A69").Value = "R_NO2_cat_(0.0, 0.005]"
